$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the "Save" header in H1, copying the style from G1 (bold/centered header style)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Fill in the new data column values
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
